# Daily Report update: 2026-02-11
# Adds the new depository snapshot (date serial 46063 = 2026-02-10) as rows
# 50-73 on Daily_Data, then refreshes the downstream Today_Summary and
# Monthly_Stats roll-ups to reflect the new figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Daily_Data: append the new day's 24 rows (12 depositories x 2 rows)
# ---------------------------------------------------------------------
$wsDaily = $wb.Worksheets.Item("Daily_Data")

$newRows = @(
    @(46063, "ASAHI DEPOSITORY LLC Registered", 24895753.652, 0, 0, 0, -942122.0600000001, 23953631.592),
    @(46063, "ASAHI DEPOSITORY LLC Eligible", 2656602.278, 0, 610906.3, -610906.3, 942122.0600000001, 2987818.038),
    @(46063, "BRINK'S, INC. Registered", 17976740.579, 0, 0, 0, -414146.13, 17562594.449),
    @(46063, "BRINK'S, INC. Eligible", 38785679.541, 0, 0, 0, 414146.13, 39199825.671),
    @(46063, "CNT DEPOSITORY, INC. Registered", 15828675.829, 0, 0, 0, -485175.71, 15343500.119),
    @(46063, "CNT DEPOSITORY, INC. Eligible", 12820940.973, 0, 368252.82, -368252.82, 485175.71, 12937863.863),
    @(46063, "DELAWARE DEPOSITORY Registered", 1966294.501, 0, 0, 0, 0, 1966294.501),
    @(46063, "DELAWARE DEPOSITORY Eligible", 15907083.262, 72286.962, 2000.6, 70286.362, 0, 15977369.624),
    @(46063, "HSBC BANK, USA Registered", 3492831.93, 0, 0, 0, -20560.25, 3472271.68),
    @(46063, "HSBC BANK, USA Eligible", 21240381.803, 0, 110629.57, -110629.57, 20560.25, 21150312.483),
    @(46063, "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered", 620749.47, 0, 0, 0, 0, 620749.47),
    @(46063, "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible", 3295246.644, 0, 0, 0, 0, 3295246.644),
    @(46063, "JP MORGAN CHASE BANK NA Registered", 12117421.82, 0, 0, 0, -82164.5, 12035257.32),
    @(46063, "JP MORGAN CHASE BANK NA Eligible", 153762189.683, 0, 2565810.5, -2565810.5, 82164.5, 151278543.683),
    @(46063, "LOOMIS INTERNATIONAL (US) LLC Registered", 7540301.607, 0, 0, 0, -166001.84, 7374299.767),
    @(46063, "LOOMIS INTERNATIONAL (US) LLC Eligible", 24268277.426, 0, 1088419.08, -1088419.08, 166001.84, 23345860.186),
    @(46063, "MALCA-AMIT ARMORED, INC. Registered", 0, 0, 0, 0, 0, 0),
    @(46063, "MALCA-AMIT ARMORED, INC. Eligible", 0, 0, 0, 0, 0, 0),
    @(46063, "MALCA-AMIT USA, LLC Registered", 1225506.264, 0, 0, 0, 0, 1225506.264),
    @(46063, "MALCA-AMIT USA, LLC Eligible", 798026.177, 0, 0, 0, 0, 798026.177),
    @(46063, "MANFRA, TORDELLA & BROOKES, LLC Registered", 8139014.54, 0, 0, 0, -1100405.34, 7038609.2),
    @(46063, "MANFRA, TORDELLA & BROOKES, LLC Eligible", 11156818.426, 0, 30489.762, -30489.762, 1100405.34, 12226734.004),
    @(46063, "STONEX PRECIOUS METALS LLC Registered", 7591598.24, 0, 0, 0, -46307.1, 7545291.14),
    @(46063, "STONEX PRECIOUS METALS LLC Eligible", 186890.28, 0, 0, 0, 46307.1, 233197.38)
)

$startRow = 50
$dateFormat = $wsDaily.Range("A2").NumberFormat

$i = 0
while ($i -lt $newRows.Count) {
    $row = $newRows[$i]
    $r = $startRow + $i

    $dateCell = $wsDaily.Cells.Item($r, 1)
    $dateCell.NumberFormat = $dateFormat
    $dateCell.Value = $row[0]

    $wsDaily.Cells.Item($r, 2).Value = $row[1]
    $wsDaily.Cells.Item($r, 3).Value = $row[2]
    $wsDaily.Cells.Item($r, 4).Value = $row[3]
    $wsDaily.Cells.Item($r, 5).Value = $row[4]
    $wsDaily.Cells.Item($r, 6).Value = $row[5]
    $wsDaily.Cells.Item($r, 7).Value = $row[6]
    $wsDaily.Cells.Item($r, 8).Value = $row[7]

    $i = $i + 1
}

# ---------------------------------------------------------------------
# 2) Today_Summary: refresh Eligible / Registered / Total_Stock per
#    depository with today's updated totals.
# ---------------------------------------------------------------------
$wsToday = $wb.Worksheets.Item("Today_Summary")

$todaySummary = @(
    @(2, 2987818.038, 23953631.592, 26941449.63),
    @(3, 39199825.671, 17562594.449, 56762420.12),
    @(4, 12937863.863, 15343500.119, 28281363.982),
    @(5, 15977369.624, 1966294.501, 17943664.125),
    @(6, 21150312.483, 3472271.68, 24622584.163),
    @(8, 151278543.683, 12035257.32, 163313801.003),
    @(9, 23345860.186, 7374299.767, 30720159.953),
    @(12, 12226734.004, 7038609.2, 19265343.204),
    @(13, 233197.38, 7545291.14, 7778488.52)
)

$i = 0
while ($i -lt $todaySummary.Count) {
    $row = $todaySummary[$i]
    $r = $row[0]
    $wsToday.Cells.Item($r, 2).Value = $row[1]
    $wsToday.Cells.Item($r, 3).Value = $row[2]
    $wsToday.Cells.Item($r, 4).Value = $row[3]
    $i = $i + 1
}

# ---------------------------------------------------------------------
# 3) Monthly_Stats: refresh the month-to-date roll-up (top block) and the
#    per-depository detail block further down the sheet.
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly_Stats")

$wsMonthly.Range("B2").Value = 283430797.753
$wsMonthly.Range("C2").Value = 98138005.502
$wsMonthly.Range("D2").Value = 381568803.255

# Detail rows: row, RECEIVED(C), WITHDRAWN(D), TOTAL_TODAY(E).
# A value of $null means that column is unchanged from before.
$monthlyDetail = @(
    @(7, $null, 1239610.7, 2987818.038),
    @(8, $null, $null, 23953631.592),
    @(9, $null, $null, 39199825.671),
    @(10, $null, $null, 17562594.449),
    @(11, $null, 2888193.293, 12937863.863),
    @(12, $null, $null, 15343500.119),
    @(13, 241770.911, 45264.589, 15977369.624),
    @(15, $null, 110629.57, 21150312.483),
    @(16, $null, $null, 3472271.68),
    @(19, $null, 5550344.199999999, 151278543.683),
    @(20, $null, $null, 12035257.32),
    @(21, $null, 2589772.75, 23345860.186),
    @(22, $null, $null, 7374299.767),
    @(27, $null, 252863.959, 12226734.004),
    @(28, $null, $null, 7038609.2),
    @(29, $null, $null, 233197.38),
    @(30, $null, $null, 7545291.14)
)

$i = 0
while ($i -lt $monthlyDetail.Count) {
    $row = $monthlyDetail[$i]
    $r = $row[0]

    if ($row[1] -ne $null) {
        $wsMonthly.Cells.Item($r, 3).Value = $row[1]
    }
    if ($row[2] -ne $null) {
        $wsMonthly.Cells.Item($r, 4).Value = $row[2]
    }
    if ($row[3] -ne $null) {
        $wsMonthly.Cells.Item($r, 5).Value = $row[3]
    }

    $i = $i + 1
}
